$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50
$ws.Range("B50").Value = 310
$ws.Range("C50").Value = -105
$ws.Range("E50").Value = -117.34

# Row 54
$ws.Range("B54").Value = 387
$ws.Range("C54").Value = 39.5
$ws.Range("E54").Value = 44.15

# Row 72
$ws.Range("B72").Value = 327
$ws.Range("C72").Value = 1154
$ws.Range("E72").Value = 1788.7

# Row 86
$ws.Range("B86").Value = 207
$ws.Range("C86").Value = 2108
$ws.Range("E86").Value = 2951.2

# Row 145
$ws.Range("B145").Value = 18
$ws.Range("C145").Value = 2
$ws.Range("E145").Value = 36.84

# Row 153
$ws.Range("B153").Value = 159
$ws.Range("C153").Value = 12.5
$ws.Range("E153").Value = 34.38

# Row 160
$ws.Range("B160").Value = 127
$ws.Range("C160").Value = 1.5
$ws.Range("E160").Value = 4.2

# Row 162
$ws.Range("B162").Value = 29
$ws.Range("C162").Value = 43.5
$ws.Range("E162").Value = 121.8

# Row 166
$ws.Range("B166").Value = 142
$ws.Range("C166").Value = 72.5
$ws.Range("E166").Value = 213.88

# Row 213
$ws.Range("B213").Value = 131
$ws.Range("C213").Value = 35
$ws.Range("E213").Value = 143.9

# Row 221
$ws.Range("B221").Value = 375
$ws.Range("C221").Value = 48
$ws.Range("E221").Value = 216

# Row 245
$ws.Range("B245").Value = 21
$ws.Range("C245").Value = 7.3
$ws.Range("E245").Value = 29.2

# Row 253
$ws.Range("B253").Value = 90
$ws.Range("C253").Value = 12
$ws.Range("E253").Value = 51

# Row 623
$ws.Range("B623").Value = 398
$ws.Range("C623").Value = 285.5
$ws.Range("E623").Value = 148.46

# Row 625
$ws.Range("B625").Value = 222
$ws.Range("C625").Value = 199.5
$ws.Range("E625").Value = 103.74

# Row 627
$ws.Range("B627").Value = 205
$ws.Range("C627").Value = 138.5
$ws.Range("E627").Value = 96.95

# Row 628
$ws.Range("B628").Value = 156
$ws.Range("C628").Value = 66.5
$ws.Range("E628").Value = 45.58

# Row 638 (Total row)
$ws.Range("C638").Value = 31058.43
$ws.Range("E638").Value = 56243.89
